$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.521.32"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "1.677.85"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5318"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2696"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.684.70"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.511"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5592"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "0.0₅8348"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "26.546.04"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.795"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.323"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +4.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.416"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06301"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.291"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.605"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.68%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.692"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.015"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6168"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.786"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.106"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.34%  "
$ws.Range("D40").Value = "1.095.81"
$ws.Range("E40").Value = "  +6.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8620"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "58.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.184"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05195"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.031"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.45%  "
